# Fruta / hortaliza, semanal
# Insert a new weekly record as row 171 (pushing the existing rows 171-177
# down to 172-178), mirroring the rest of the row's descriptive columns
# from the entry that follows it, but with its own measurement values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new record at row 171; everything currently at/after
# row 171 (old rows 171-177) shifts down to 172-178.
$ws.Rows.Item(171).Insert()

$ws.Range("A171").Value = 1
$ws.Range("B171").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C171").Value = "Arica y Parinacota"
$ws.Range("D171").Value = 44855
$ws.Range("E171").Value = 15
$ws.Range("F171").Value = "Fruta"
$ws.Range("G171").Value = 100106
$ws.Range("H171").Value = "Oleaginosos"
$ws.Range("I171").Value = 100106002
$ws.Range("J171").Value = "Palta"
$ws.Range("K171").Value = "Fuerte"
$ws.Range("L171").Value = "Segunda"
$ws.Range("M171").Value = 600
$ws.Range("N171").Value = 9000
$ws.Range("O171").Value = 10000
$ws.Range("P171").Value = 9583
$ws.Range("Q171").Value = "$/bandeja 10 kilos"
$ws.Range("R171").Value = "Perú"
$ws.Range("S171").Value = 958
$ws.Range("T171").Value = 10
